$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (date moved from 08-04 to 08-05)
$ws.Name = "Through 2022-08-05"

# Update the "August (through 08-04)" label to "August (through 08-05)"
$ws.Range("A9").Value = "August (through 08-05)"

# Update August row (row 9) values for each year column (B..I)
$ws.Range("B9").Value = 6
$ws.Range("C9").Value = 6
$ws.Range("D9").Value = 11
$ws.Range("E9").Value = 12
$ws.Range("F9").Value = 8
$ws.Range("G9").Value = 30
$ws.Range("H9").Value = 25
$ws.Range("I9").Value = 25

# Update Total row (row 10) values for each year column (B..I)
$ws.Range("B10").Value = 168
$ws.Range("C10").Value = 308
$ws.Range("D10").Value = 476
$ws.Range("E10").Value = 437
$ws.Range("F10").Value = 312
$ws.Range("G10").Value = 651
$ws.Range("H10").Value = 935
$ws.Range("I10").Value = 995
